# Auto-generated edit script: updates market-price / profit columns (H:N)
# on the Tonberry Profits leve-crafting sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 1893.5834
$ws.Range("I6").Value = 80.75
$ws.Range("J6").Value = 2800
$ws.Range("K6").Value = 242.25
$ws.Range("L6").Value = 8400
$ws.Range("M6").Value = -130.25
$ws.Range("N6").Value = -8624
# Row 18
$ws.Range("H18").Value = 15405.5625
$ws.Range("I18").Value = 13642.857
$ws.Range("K18").Value = 13642.857
$ws.Range("M18").Value = -13358.857
# Row 32
$ws.Range("H32").Value = 766.6667
$ws.Range("J32").Value = 967
$ws.Range("L32").Value = 967
$ws.Range("N32").Value = -1619
# Row 70
$ws.Range("H70").Value = 38316.668
$ws.Range("J70").Value = 37250
$ws.Range("L70").Value = 111750
$ws.Range("N70").Value = -112290
# Row 73
$ws.Range("H73").Value = 38316.668
$ws.Range("J73").Value = 37250
$ws.Range("L73").Value = 111750
$ws.Range("N73").Value = -113622
# Row 86
$ws.Range("H86").Value = 2555.7144
$ws.Range("I86").Value = 2400
$ws.Range("K86").Value = 2400
$ws.Range("M86").Value = -1277
# Row 89
$ws.Range("H89").Value = 2555.7144
$ws.Range("I89").Value = 2400
$ws.Range("K89").Value = 12000
$ws.Range("M89").Value = -6384
# Row 125
$ws.Range("H125").Value = 874
$ws.Range("J125").Value = 967.5
$ws.Range("L125").Value = 8707.5
$ws.Range("N125").Value = -13627.5
# Row 135
$ws.Range("H135").Value = 694.1111
$ws.Range("I135").Value = 593.375
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 5340.375
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -2805.375
$ws.Range("N135").Value = -18570
# Row 139
$ws.Range("H139").Value = 45000
$ws.Range("J139").Value = 45000
$ws.Range("L139").Value = 45000
$ws.Range("N139").Value = -55280

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3693.0535
$ws.Range("I32").Value = 2279
$ws.Range("J32").Value = 8877.916999999999
$ws.Range("K32").Value = 2279
$ws.Range("L32").Value = 8877.916999999999
$ws.Range("M32").Value = -1992
$ws.Range("N32").Value = -9451.916999999999
# Row 88
$ws.Range("H88").Value = 3193.9
$ws.Range("I88").Value = 2168
$ws.Range("J88").Value = 4219.8
$ws.Range("K88").Value = 2168
$ws.Range("L88").Value = 4219.8
$ws.Range("M88").Value = -1762
$ws.Range("N88").Value = -5031.8
# Row 91
$ws.Range("H91").Value = 3193.9
$ws.Range("I91").Value = 2168
$ws.Range("J91").Value = 4219.8
$ws.Range("K91").Value = 2168
$ws.Range("L91").Value = 4219.8
$ws.Range("M91").Value = -764
$ws.Range("N91").Value = -7027.8
# Row 110
$ws.Range("H110").Value = 1790.125
$ws.Range("I110").Value = 1720.1666
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1720.1666
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 324.8334
$ws.Range("N110").Value = -6090
# Row 122
$ws.Range("H122").Value = 1529.7428
$ws.Range("I122").Value = 1488.8667
$ws.Range("K122").Value = 4466.6001
$ws.Range("M122").Value = -2016.6001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 2312
$ws.Range("J15").Value = 2312
$ws.Range("L15").Value = 2312
$ws.Range("N15").Value = -2652
# Row 31
$ws.Range("H31").Value = 2386.1667
$ws.Range("I31").Value = 891.36
$ws.Range("J31").Value = 5783.4546
$ws.Range("K31").Value = 891.36
$ws.Range("L31").Value = 5783.4546
$ws.Range("M31").Value = -596.36
$ws.Range("N31").Value = -6373.4546
# Row 34
$ws.Range("H34").Value = 2386.1667
$ws.Range("I34").Value = 891.36
$ws.Range("J34").Value = 5783.4546
$ws.Range("K34").Value = 891.36
$ws.Range("L34").Value = 5783.4546
$ws.Range("M34").Value = -689.36
$ws.Range("N34").Value = -6187.4546
# Row 99
$ws.Range("H99").Value = 1871.125
$ws.Range("I99").Value = 1495.6666
$ws.Range("J99").Value = 2997.5
$ws.Range("K99").Value = 1495.6666
$ws.Range("L99").Value = 2997.5
$ws.Range("M99").Value = 2.333399999999983
$ws.Range("N99").Value = -5993.5
# Row 107
$ws.Range("H107").Value = 534.55554
$ws.Range("I107").Value = 443.0909
$ws.Range("J107").Value = 678.2857
$ws.Range("K107").Value = 443.0909
$ws.Range("L107").Value = 678.2857
$ws.Range("M107").Value = 1476.9091
$ws.Range("N107").Value = -4518.2857
# Row 126
$ws.Range("H126").Value = 1871.125
$ws.Range("I126").Value = 1495.6666
$ws.Range("J126").Value = 2997.5
$ws.Range("K126").Value = 4486.9998
$ws.Range("L126").Value = 8992.5
$ws.Range("M126").Value = -2016.9998
$ws.Range("N126").Value = -13932.5
# Row 134
$ws.Range("H134").Value = 1112
$ws.Range("I134").Value = 1093.3
$ws.Range("K134").Value = 3279.9
$ws.Range("M134").Value = -744.8999999999996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 591.4167
$ws.Range("I5").Value = 519.2
$ws.Range("K5").Value = 1557.6
$ws.Range("M5").Value = -1445.6
# Row 7
$ws.Range("H7").Value = 689.8
$ws.Range("I7").Value = 49
$ws.Range("K7").Value = 147
$ws.Range("M7").Value = -35
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
# Row 92
$ws.Range("H92").Value = 433.55554
$ws.Range("J92").Value = 456.25
$ws.Range("L92").Value = 1368.75
$ws.Range("N92").Value = -3864.75
# Row 107
$ws.Range("H107").Value = 594.375
$ws.Range("J107").Value = 594.375
$ws.Range("L107").Value = 1783.125
$ws.Range("N107").Value = -5623.125
# Row 122
$ws.Range("H122").Value = 1048.8334
$ws.Range("J122").Value = 1896.25
$ws.Range("L122").Value = 17066.25
$ws.Range("N122").Value = -21966.25
# Row 131
$ws.Range("H131").Value = 13179599
$ws.Range("J131").Value = 27346.334
$ws.Range("L131").Value = 82039.00199999999
$ws.Range("N131").Value = -92119.00199999999
# Row 135
$ws.Range("H135").Value = 591.4167
$ws.Range("I135").Value = 519.2
$ws.Range("K135").Value = 4672.8
$ws.Range("M135").Value = -2137.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 18425
$ws.Range("I46").Value = 13900
$ws.Range("J46").Value = 19933.334
$ws.Range("K46").Value = 13900
$ws.Range("L46").Value = 19933.334
$ws.Range("M46").Value = -13744
$ws.Range("N46").Value = -20245.334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1441.125
$ws.Range("J22").Value = 1450.6923
$ws.Range("L22").Value = 1450.6923
$ws.Range("N22").Value = -2040.6923
# Row 27
$ws.Range("H27").Value = 1441.125
$ws.Range("J27").Value = 1450.6923
$ws.Range("L27").Value = 1450.6923
$ws.Range("N27").Value = -1664.6923
# Row 46
$ws.Range("H46").Value = 1936.2727
$ws.Range("I46").Value = 1299
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1299
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1111
$ws.Range("N46").Value = -2376
# Row 82
$ws.Range("H82").Value = 3439.0833
$ws.Range("I82").Value = 1947.25
$ws.Range("J82").Value = 4185
$ws.Range("K82").Value = 1947.25
$ws.Range("L82").Value = 4185
$ws.Range("M82").Value = -1586.25
$ws.Range("N82").Value = -4907
# Row 85
$ws.Range("H85").Value = 3439.0833
$ws.Range("I85").Value = 1947.25
$ws.Range("J85").Value = 4185
$ws.Range("K85").Value = 1947.25
$ws.Range("L85").Value = 4185
$ws.Range("M85").Value = -699.25
$ws.Range("N85").Value = -6681
# Row 93
$ws.Range("H93").Value = 1994
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1994
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1994
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -4490

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 16999.4
$ws.Range("J31").Value = 16999.4
$ws.Range("L31").Value = 16999.4
$ws.Range("N31").Value = -17695.4
# Row 49
$ws.Range("H49").Value = 42178.285
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
# Row 62
$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 2750
$ws.Range("K62").Value = 2750
$ws.Range("M62").Value = -2126
# Row 65
$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 2750
$ws.Range("K65").Value = 13750
$ws.Range("M65").Value = -10630
# Row 70
$ws.Range("H70").Value = 30000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
# Row 73
$ws.Range("H73").Value = 30000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
# Row 107
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 2400
$ws.Range("M107").Value = -480
# Row 136
$ws.Range("H136").Value = 3623.7273
$ws.Range("I136").Value = 3602.6667
$ws.Range("J136").Value = 3649
$ws.Range("K136").Value = 10808.0001
$ws.Range("L136").Value = 10947
$ws.Range("M136").Value = -8258.000100000001
$ws.Range("N136").Value = -16047

